$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (shifts G:O right to H:P)
$ws.Range("F1").EntireColumn.Insert()

# Fill in new column F values
$ws.Range("F2").Value = "Interactions"
$ws.Range("F3").Value = "None"
$ws.Range("F6").Value = "None"
$ws.Range("F9").Value = "None"
$ws.Range("F13").Value = "None"

# Match the saved selection / view state
$ws.Range("F14").Select()
$excel.ActiveWindow.Zoom = 130
